# Update NATMI TPM-derived Cxcl13-Cxcr5 LR-pair stats.
# The underlying TPM recomputation changes every specificity/weight
# column for the existing FAPs/MuSCs/Resolving-Mac x (ECs/FAPs/MuSCs/
# Resolving-Mac) pairs, reorders a couple of target-cluster rows, and
# appends the three Resolving-Mac sending-cluster rows that were missing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 9).Value = 0.5162107379131895
$ws.Cells.Item(2, 10).Value = 0.5162107379131895
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.2320226666666667
$ws.Cells.Item(2, 14).Value = 0.696068
$ws.Cells.Item(2, 15).Value = 0.2267671859540574
$ws.Cells.Item(2, 16).Value = 0.2267671859540574
$ws.Cells.Item(2, 17).Value = 0.9388716772142222
$ws.Cells.Item(2, 18).Value = 8.449845094928
$ws.Cells.Item(2, 19).Value = 0.1170596563958414
$ws.Cells.Item(2, 20).Value = 0.1170596563958414

# Row 3
$ws.Cells.Item(3, 9).Value = 0.5162107379131895
$ws.Cells.Item(3, 10).Value = 0.5162107379131895
$ws.Cells.Item(3, 13).Value = 0.5973453333333333
$ws.Cells.Item(3, 15).Value = 0.5838150307848733
$ws.Cells.Item(3, 16).Value = 0.5838150307848734
$ws.Cells.Item(3, 17).Value = 2.417137183361777
$ws.Cells.Item(3, 19).Value = 0.3013715878462709
$ws.Cells.Item(3, 20).Value = 0.3013715878462709

# Row 4
$ws.Cells.Item(4, 9).Value = 0.5162107379131895
$ws.Cells.Item(4, 10).Value = 0.5162107379131895
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.120982
$ws.Cells.Item(4, 14).Value = 0.362946
$ws.Cells.Item(4, 15).Value = 0.1182416704593248
$ws.Cells.Item(4, 16).Value = 0.1182416704593248
$ws.Cells.Item(4, 17).Value = 0.4895494689573333
$ws.Cells.Item(4, 18).Value = 4.405945220616
$ws.Cells.Item(4, 19).Value = 0.06103761995989625
$ws.Cells.Item(4, 20).Value = 0.06103761995989625

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.046465333333333
$ws.Cells.Item(5, 8).Value = 12.139396
$ws.Cells.Item(5, 9).Value = 0.5162107379131895
$ws.Cells.Item(5, 10).Value = 0.5162107379131895
$ws.Cells.Item(5, 13).Value = 0.07282566666666666
$ws.Cells.Item(5, 14).Value = 0.218477
$ws.Cells.Item(5, 15).Value = 0.07117611280174437
$ws.Cells.Item(5, 16).Value = 0.07117611280174438
$ws.Cells.Item(5, 17).Value = 0.2946865355435555
$ws.Cells.Item(5, 18).Value = 2.652178819892
$ws.Cells.Item(5, 19).Value = 0.03674187371118087
$ws.Cells.Item(5, 20).Value = 0.03674187371118088

# Row 6
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 7).Value = 0.155986
$ws.Cells.Item(6, 8).Value = 0.467958
$ws.Cells.Item(6, 9).Value = 0.01989925565426652
$ws.Cells.Item(6, 10).Value = 0.01989925565426652
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.2320226666666667
$ws.Cells.Item(6, 14).Value = 0.696068
$ws.Cells.Item(6, 15).Value = 0.2267671859540574
$ws.Cells.Item(6, 16).Value = 0.2267671859540574
$ws.Cells.Item(6, 17).Value = 0.03619228768266666
$ws.Cells.Item(6, 18).Value = 0.325730589144
$ws.Cells.Item(6, 19).Value = 0.004512498207298384
$ws.Cells.Item(6, 20).Value = 0.004512498207298384

# Row 7
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 7).Value = 0.155986
$ws.Cells.Item(7, 8).Value = 0.467958
$ws.Cells.Item(7, 9).Value = 0.01989925565426652
$ws.Cells.Item(7, 10).Value = 0.01989925565426652
$ws.Cells.Item(7, 13).Value = 0.5973453333333333
$ws.Cells.Item(7, 14).Value = 1.792036
$ws.Cells.Item(7, 15).Value = 0.5838150307848733
$ws.Cells.Item(7, 16).Value = 0.5838150307848734
$ws.Cells.Item(7, 17).Value = 0.09317750916533332
$ws.Cells.Item(7, 18).Value = 0.838597582488
$ws.Cells.Item(7, 19).Value = 0.01161748455239167
$ws.Cells.Item(7, 20).Value = 0.01161748455239167

# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 4).Value = 'MuSCs'
$ws.Cells.Item(8, 7).Value = 0.155986
$ws.Cells.Item(8, 8).Value = 0.467958
$ws.Cells.Item(8, 9).Value = 0.01989925565426652
$ws.Cells.Item(8, 10).Value = 0.01989925565426652
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.120982
$ws.Cells.Item(8, 14).Value = 0.362946
$ws.Cells.Item(8, 15).Value = 0.1182416704593248
$ws.Cells.Item(8, 16).Value = 0.1182416704593248
$ws.Cells.Item(8, 17).Value = 0.018871498252
$ws.Cells.Item(8, 18).Value = 0.169843484268
$ws.Cells.Item(8, 19).Value = 0.002352921229457638
$ws.Cells.Item(8, 20).Value = 0.002352921229457638

# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 7).Value = 0.155986
$ws.Cells.Item(9, 8).Value = 0.467958
$ws.Cells.Item(9, 9).Value = 0.01989925565426652
$ws.Cells.Item(9, 10).Value = 0.01989925565426652
$ws.Cells.Item(9, 13).Value = 0.07282566666666666
$ws.Cells.Item(9, 14).Value = 0.218477
$ws.Cells.Item(9, 15).Value = 0.07117611280174437
$ws.Cells.Item(9, 16).Value = 0.07117611280174438
$ws.Cells.Item(9, 17).Value = 0.01135978444066667
$ws.Cells.Item(9, 18).Value = 0.102238059966
$ws.Cells.Item(9, 19).Value = 0.001416351665118823
$ws.Cells.Item(9, 20).Value = 0.001416351665118823

# Row 10
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.636334333333334
$ws.Cells.Item(10, 8).Value = 10.909003
$ws.Cells.Item(10, 9).Value = 0.463890006432544
$ws.Cells.Item(10, 10).Value = 0.463890006432544
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.2320226666666667
$ws.Cells.Item(10, 14).Value = 0.696068
$ws.Cells.Item(10, 15).Value = 0.2267671859540574
$ws.Cells.Item(10, 16).Value = 0.2267671859540574
$ws.Cells.Item(10, 17).Value = 0.8437119889115556
$ws.Cells.Item(10, 18).Value = 7.593407900204
$ws.Cells.Item(10, 19).Value = 0.1051950313509176
$ws.Cells.Item(10, 20).Value = 0.1051950313509176

# Row 11
$ws.Cells.Item(11, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 2).Value = 'Cxcl13'
$ws.Cells.Item(11, 3).Value = 'Cxcr5'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.636334333333334
$ws.Cells.Item(11, 8).Value = 10.909003
$ws.Cells.Item(11, 9).Value = 0.463890006432544
$ws.Cells.Item(11, 10).Value = 0.463890006432544
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.5973453333333333
$ws.Cells.Item(11, 14).Value = 1.792036
$ws.Cells.Item(11, 15).Value = 0.5838150307848733
$ws.Cells.Item(11, 16).Value = 0.5838150307848734
$ws.Cells.Item(11, 17).Value = 2.172147344456445
$ws.Cells.Item(11, 18).Value = 19.549326100108
$ws.Cells.Item(11, 19).Value = 0.2708259583862108
$ws.Cells.Item(11, 20).Value = 0.2708259583862108

# Row 12
$ws.Cells.Item(12, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(12, 2).Value = 'Cxcl13'
$ws.Cells.Item(12, 3).Value = 'Cxcr5'
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.636334333333334
$ws.Cells.Item(12, 8).Value = 10.909003
$ws.Cells.Item(12, 9).Value = 0.463890006432544
$ws.Cells.Item(12, 10).Value = 0.463890006432544
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.120982
$ws.Cells.Item(12, 14).Value = 0.362946
$ws.Cells.Item(12, 15).Value = 0.1182416704593248
$ws.Cells.Item(12, 16).Value = 0.1182416704593248
$ws.Cells.Item(12, 17).Value = 0.4399310003153333
$ws.Cells.Item(12, 18).Value = 3.959379002838
$ws.Cells.Item(12, 19).Value = 0.05485112926997095
$ws.Cells.Item(12, 20).Value = 0.05485112926997094

# Row 13
$ws.Cells.Item(13, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 2).Value = 'Cxcl13'
$ws.Cells.Item(13, 3).Value = 'Cxcr5'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.636334333333334
$ws.Cells.Item(13, 8).Value = 10.909003
$ws.Cells.Item(13, 9).Value = 0.463890006432544
$ws.Cells.Item(13, 10).Value = 0.463890006432544
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.07282566666666666
$ws.Cells.Item(13, 14).Value = 0.218477
$ws.Cells.Item(13, 15).Value = 0.07117611280174437
$ws.Cells.Item(13, 16).Value = 0.07117611280174438
$ws.Cells.Item(13, 17).Value = 0.2648184720478889
$ws.Cells.Item(13, 18).Value = 2.383366248431
$ws.Cells.Item(13, 19).Value = 0.03301788742544468
$ws.Cells.Item(13, 20).Value = 0.03301788742544468

